# Commit: "finish the main file"
# The Answers column (C) originally mixed "ct" and "mri" values.
# All remaining "ct" answers are corrected to "mri" so every row in
# the VQA answer sheet now reads "mri".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq "ct") {
        $cell.Value2 = "mri"
    }
}
